# poster v4.1 (final final version)
# Remove the "Picture 39" image (the small cartoon-face / Smiling Globe icon
# placed next to the poster title) from the one and only slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Picture 39" -and $shp.Id -eq 40) {
        $shp.Delete()
        break
    }
}
